# Insert a new row at position 41, shifting existing data rows 41-153 down to 42-154.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with the new weekly data point.
$ws.Range("A41").Value2 = 11
$ws.Range("B41").Value2 = "Vega Monumental Concepción"
$ws.Range("C41").Value2 = "Bíobío"
$ws.Range("D41").Value2 = 45260
$ws.Range("E41").Value2 = 8
$ws.Range("F41").Value2 = 100112037
$ws.Range("G41").Value2 = "Cebollín"
$ws.Range("H41").Value2 = "Sin especificar"
$ws.Range("I41").Value2 = "Primera"
$ws.Range("J41").Value2 = 80
$ws.Range("K41").Value2 = 3500
$ws.Range("L41").Value2 = 3500
$ws.Range("M41").Value2 = 3500
$ws.Range("N41").Value2 = "`$/paquete 36 unidades"
$ws.Range("O41").Value2 = "Región Metropolitana"
$ws.Range("P41").Value2 = 97
$ws.Range("Q41").Value2 = 36
$ws.Range("R41").Value2 = "Hortaliza"
